# New Project to Push
# Adds a "PRD" row to the URL sheet and a second/third User ID+Password
# block (SteveH/Tech@01, JohnR/Mech@03) to the ID_PWD sheet, then leaves
# the ID_PWD sheet as the active tab/view.

$wb = $excel.ActiveWorkbook
$wsUrl = $wb.Worksheets.Item("URL")
$wsId  = $wb.Worksheets.Item("ID_PWD")

# ---------------------------------------------------------------------
# Sheet "URL": append the PRD row (row 4)
# ---------------------------------------------------------------------
$wsUrl.Cells.Item(4, 1).Value = "PRD"

$prdUrl = "https://prd.mobolutions.com:1443/sap/bc/ui2/flp#Shell-home"
$wsUrl.Cells.Item(4, 2).Value = $prdUrl
$wsUrl.Hyperlinks.Add($wsUrl.Cells.Item(4, 2), "https://prd.mobolutions.com:1443/sap/bc/ui2/flp", "Shell-home", "", $prdUrl) | Out-Null

# Re-apply the same look as the other URL cells (B2/B3) without creating
# a brand-new duplicate cell style: copy formatting only from B2.
$wsUrl.Cells.Item(2, 2).Copy()
$wsUrl.Cells.Item(4, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "ID_PWD": add two more User ID / Password column pairs
# (D:E = SteveH/Tech@01, F:G = JohnR/Mech@03)
# ---------------------------------------------------------------------

# Headers row 1 (copy header formatting from the existing B1/C1 headers)
$wsId.Cells.Item(1, 4).Value = "User ID"
$wsId.Cells.Item(1, 2).Copy()
$wsId.Cells.Item(1, 4).PasteSpecial(-4122)

$wsId.Cells.Item(1, 5).Value = "Password"
$wsId.Cells.Item(1, 3).Copy()
$wsId.Cells.Item(1, 5).PasteSpecial(-4122)

$wsId.Cells.Item(1, 6).Value = "User ID"
$wsId.Cells.Item(1, 2).Copy()
$wsId.Cells.Item(1, 6).PasteSpecial(-4122)

$wsId.Cells.Item(1, 7).Value = "Password"
$wsId.Cells.Item(1, 3).Copy()
$wsId.Cells.Item(1, 7).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row 2
$wsId.Cells.Item(2, 4).Value = "SteveH"
$wsId.Cells.Item(2, 6).Value = "JohnR"

$wsId.Cells.Item(2, 5).Value = "Tech@01"
$wsId.Hyperlinks.Add($wsId.Cells.Item(2, 5), "mailto:Tech@01") | Out-Null
$wsId.Cells.Item(2, 3).Copy()
$wsId.Cells.Item(2, 5).PasteSpecial(-4122)

$wsId.Cells.Item(2, 7).Value = "Mech@03"
$wsId.Hyperlinks.Add($wsId.Cells.Item(2, 7), "mailto:Mech@03") | Out-Null
$wsId.Cells.Item(2, 3).Copy()
$wsId.Cells.Item(2, 7).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Best-fit the new Password column (G) like the other data columns.
$wsId.Columns.Item(7).EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# View state: ID_PWD becomes the selected/active sheet & tab.
# ---------------------------------------------------------------------
$wsUrl.Activate()
$wsUrl.Range("E7").Select()

$wsId.Activate()
$wsId.Range("K14").Select()
